$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 192.352885
$ws.Range("H2").Value = 577.058655
$ws.Range("I2").Value = 0.3381139681403721
$ws.Range("J2").Value = 0.3381139681403722
$ws.Range("M2").Value = 211.980367
$ws.Range("N2").Value = 635.9411009999999
$ws.Range("O2").Value = 0.9885149156420702
$ws.Range("P2").Value = 0.9885149156420702
$ws.Range("Q2").Value = 40775.03515580879
$ws.Range("R2").Value = 366975.3164022791
$ws.Range("S2").Value = 0.3342307006936855
$ws.Range("T2").Value = 0.3342307006936856

# Row 3
$ws.Range("G3").Value = 192.352885
$ws.Range("H3").Value = 577.058655
$ws.Range("I3").Value = 0.3381139681403721
$ws.Range("J3").Value = 0.3381139681403722
$ws.Range("O3").Value = 0.003992992409159323
$ws.Range("P3").Value = 0.003992992409159324
$ws.Range("Q3").Value = 164.7060689565783
$ws.Range("R3").Value = 1482.354620609205
$ws.Range("S3").Value = 0.001350086508215243
$ws.Range("T3").Value = 0.001350086508215243

# Row 4
$ws.Range("G4").Value = 192.352885
$ws.Range("H4").Value = 577.058655
$ws.Range("I4").Value = 0.3381139681403721
$ws.Range("J4").Value = 0.3381139681403722
$ws.Range("O4").Value = 0.007492091948770576
$ws.Range("P4").Value = 0.007492091948770576
$ws.Range("Q4").Value = 309.0396591570367
$ws.Range("R4").Value = 2781.35693241333
$ws.Range("S4").Value = 0.002533180938471353
$ws.Range("T4").Value = 0.002533180938471353

# Row 5
$ws.Range("I5").Value = 0.2958833255212922
$ws.Range("J5").Value = 0.2958833255212922
$ws.Range("M5").Value = 211.980367
$ws.Range("N5").Value = 635.9411009999999
$ws.Range("O5").Value = 0.9885149156420702
$ws.Range("P5").Value = 0.9885149156420702
$ws.Range("Q5").Value = 35682.20818117612
$ws.Range("R5").Value = 321139.873630585
$ws.Range("S5").Value = 0.2924850805675753
$ws.Range("T5").Value = 0.2924850805675753

# Row 6
$ws.Range("I6").Value = 0.2958833255212922
$ws.Range("J6").Value = 0.2958833255212922
$ws.Range("O6").Value = 0.003992992409159323
$ws.Range("P6").Value = 0.003992992409159324
$ws.Range("S6").Value = 0.001181459872803337
$ws.Range("T6").Value = 0.001181459872803337

# Row 7
$ws.Range("I7").Value = 0.2958833255212922
$ws.Range("J7").Value = 0.2958833255212922
$ws.Range("O7").Value = 0.007492091948770576
$ws.Range("P7").Value = 0.007492091948770576
$ws.Range("S7").Value = 0.002216785080913537
$ws.Range("T7").Value = 0.002216785080913537

# Row 8
$ws.Range("I8").Value = 0.3660027063383355
$ws.Range("J8").Value = 0.3660027063383356
$ws.Range("M8").Value = 211.980367
$ws.Range("N8").Value = 635.9411009999999
$ws.Range("O8").Value = 0.9885149156420702
$ws.Range("P8").Value = 0.9885149156420702
$ws.Range("Q8").Value = 44138.29248211067
$ws.Range("R8").Value = 397244.632338996
$ws.Range("S8").Value = 0.3617991343808091
$ws.Range("T8").Value = 0.3617991343808092

# Row 9
$ws.Range("I9").Value = 0.3660027063383355
$ws.Range("J9").Value = 0.3660027063383356
$ws.Range("O9").Value = 0.003992992409159323
$ws.Range("P9").Value = 0.003992992409159324
$ws.Range("S9").Value = 0.001461446028140742
$ws.Range("T9").Value = 0.001461446028140743

# Row 10
$ws.Range("I10").Value = 0.3660027063383355
$ws.Range("J10").Value = 0.3660027063383356
$ws.Range("O10").Value = 0.007492091948770576
$ws.Range("P10").Value = 0.007492091948770576
$ws.Range("S10").Value = 0.002742125929385685
$ws.Range("T10").Value = 0.002742125929385686

